$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the existing date-formatted cell (A120) as a style template for the new date cells,
# so we reuse the existing style index rather than create a brand new one.
$dateTemplate = $ws.Range("A120")

# Row 121: 9/28/2021 sample, using existing remark "CRM opened 9/24/2021"
$dateTemplate.Copy()
$ws.Range("A121").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(121, 1).Value = 44467
$ws.Cells.Item(121, 2).Value = 2250.79045167291
$ws.Cells.Item(121, 3).Value = 2230.52
$ws.Range("D121").Formula = "=100*(B121-C121)/C121"
$ws.Cells.Item(121, 5).Value = 183
$ws.Cells.Item(121, 6).Value = "CRM opened 9/24/2021"

# Row 122: 9/29/2021 sample, using new remark "CRM opened 9/29/2021"
$dateTemplate.Copy()
$ws.Range("A122").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(122, 1).Value = 44468
$ws.Cells.Item(122, 2).Value = 2249.1843699999999
$ws.Cells.Item(122, 3).Value = 2230.52
$ws.Range("D122").Formula = "=100*(B122-C122)/C122"
$ws.Cells.Item(122, 5).Value = 191
$ws.Cells.Item(122, 6).Value = "CRM opened 9/29/2021"

$excel.CutCopyMode = 0

# Update the frozen pane / selection view to match the new bottom of data
$ws.Activate()
$ws.Range("A122").Select()
$excel.ActiveWindow.ScrollRow = 110
